$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report drops the oldest period column (old column D, "6 ماهه منتهی به
# 1399/06") and appends a new most-recent period column at the end
# ("12 ماهه منتهی به 1401/12"), i.e. the whole quarterly cumulative table is
# rolled forward by one column. Deleting column D shifts every other period
# column (old E..M) one slot to the left (new D..L) together with its data,
# formatting and column widths.
$ws.Columns.Item(4).Delete()

# Re-create the 13th (M) column that the delete left behind by cloning the
# formatting of the new last data column (L) into it.
$ws.Range("L1:L28").Copy()
$ws.Range("M1:M28").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New column width for M, matching the recurring 28/29 width pattern.
$ws.Columns.Item(13).ColumnWidth = 28.17

# Period header (row 8) and publish-date header (row 9) for the new column.
$ws.Cells.Item(8, 13).Value = "12 ماهه منتهی به 1401/12"
$ws.Cells.Item(9, 13).Value = "1402-02-30 (2)"

# New cumulative figures for the new "12 ماهه منتهی به 1401/12" period.
$ws.Cells.Item(11, 13).Value = 12813639
$ws.Cells.Item(12, 13).Value = -8536412
$ws.Cells.Item(13, 13).Value = 4277227
$ws.Cells.Item(14, 13).Value = -892802
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(16, 13).Value = -8591
$ws.Cells.Item(17, 13).Value = 3375834
$ws.Cells.Item(18, 13).Value = -108052
$ws.Cells.Item(19, 13).Value = 48612
$ws.Cells.Item(20, 13).Value = 3316394
$ws.Cells.Item(21, 13).Value = -398868
$ws.Cells.Item(22, 13).Value = 2917526
$ws.Cells.Item(23, 13).Value = 0
$ws.Cells.Item(24, 13).Value = 2917526
$ws.Cells.Item(25, 13).Value = 651
$ws.Cells.Item(26, 13).Value = 4484000
$ws.Cells.Item(27, 13).Value = 651
